$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B (ASIN). This shifts ASIN and all
# subsequent columns one position to the right.
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week start dates for W1..W16 (2025-01-05, incrementing by 7 days each week)
$weekStartDates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    # Update Week value in column A, removing the leading zero (W01 -> W1)
    $ws.Cells.Item($row, 1).Value = "W" + ($i + 1)

    # Force column B to be stored as text so the date-like string isn't
    # auto-converted into a date serial number.
    $cell = $ws.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $weekStartDates[$i]
}
